$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update scenario name and nominal transfer value on row 6 first, so the
# revised shared-string table keeps this text in the slot formerly used by
# "Nominal Transfer > max".
$ws.Range("A6").Value = "Saldo Kurang"
$ws.Range("D6").Value = 5000000

# Re-apply the number formatting that the value assignment above reset,
# by pulling it back from a sibling cell that already carries it.
$ws.Range("D7").Copy()
$ws.Range("D6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update header label for the Scenario column (new shared string, appended
# at the end of the table).
$ws.Range("A1").Value = "Scenario BS+ ke BS"

# Move active selection to D14 as reflected in the saved view state
$ws.Range("D14").Select()
